# ETF.xlsx restructuring:
#   - Rotate the three tabs: old "Bond" (pos 1) becomes "Equity",
#     old "Alternative" (pos 2) becomes "Bond", old "Equity" (pos 3)
#     becomes "Alternative" -- i.e. each physical sheet's data is replaced
#     by the data that belongs under its new tab name (with a handful of
#     corrected values along the way), and the tabs are renamed to match.

$wb = $excel.ActiveWorkbook

$sheetEquity = $wb.Worksheets.Item(1)   # was "Bond"      -> becomes "Equity"
$sheetBond   = $wb.Worksheets.Item(2)   # was "Alternative" -> becomes "Bond"
$sheetAlt    = $wb.Worksheets.Item(3)   # was "Equity"    -> becomes "Alternative"

# ---------------------------------------------------------------------
# New "Equity" sheet (physical sheet 1) -- 7 data rows (rows 2-8)
# ---------------------------------------------------------------------
$equityData = @(
    ,@("VOO","Traditional Equity","Vanguard 500 Index Fund","NYSEArca","Large Blend",0.0004,0.0162,18.7)
    ,@("FLCA","Traditional Equity","Franklin FTSE Canada ETF","NYSEArca","Miscellaneous Region",0.0009,0.0317,20.54)
    ,@("FLJP","Traditional Equity","Franklin FTSE Japan ETF","NYSEArca","Japan Stock",0.0009,0.0224,15.34)
    ,@("FLAU","Traditional Equity","Franklin FTSE Australia ETF","NYSEArca","Miscellaneous Region",0.0009,0.0472,23.72)
    ,@("FLKR","Traditional Equity","Franklin FTSE South Korea ETF","NYSEArca","Miscellaneous Region",0.0009,0.0375,25.41)
    ,@("SPEU","Traditional Equity","SPDR Portfolio Europe ETF","NYSEArca","Europe Stock",0.0007,0.0333,19.53)
    ,@("SPEM","Traditional Equity","SPDR Portfolio Emerging Markets ETF","NYSEArca","Diversified Emerging Mkts",0.0007,0.0339,17.73)
)

for ($i = 0; $i -lt $equityData.Count; $i++) {
    $row = $i + 2
    $rec = $equityData[$i]
    for ($c = 0; $c -lt 8; $c++) {
        $sheetEquity.Cells.Item($row, $c + 1).Value = $rec[$c]
    }
}

# ---------------------------------------------------------------------
# New "Bond" sheet (physical sheet 2) -- 4 data rows (rows 2-5)
# ---------------------------------------------------------------------
$bondData = @(
    ,@("FLIA","Traditional Bond","Franklin International Aggregate Bond ETF","Cboe US","",0.0025,0.1517,4.26)
    ,@("BILS","Traditional Bond","SPDR Bloomberg 3-12 Month T-Bill ETF","NYSEArca","Ultrashort Bond",0.00135,0.0457,0)
    ,@("VCLT","Traditional Bond","Vanguard Long-Term Corporate Bond Index Fund","NasdaqGM","Long-Term Bond",0.0007,0.0535,13.98)
    ,@("VWOB","Traditional Bond","Vanguard Emerging Markets Government Bond Index Fund","NasdaqGM","Emerging Markets Bond",0.002,0.0597,11.02)
)

for ($i = 0; $i -lt $bondData.Count; $i++) {
    $row = $i + 2
    $rec = $bondData[$i]
    for ($c = 0; $c -lt 8; $c++) {
        $sheetBond.Cells.Item($row, $c + 1).Value = $rec[$c]
    }
}

# ---------------------------------------------------------------------
# New "Alternative" sheet (physical sheet 3) -- 4 data rows (rows 2-5);
# shrinks from 8 rows down to 5, so clear the now-unused rows 6-8 first.
# ---------------------------------------------------------------------
$sheetAlt.Range("A6:J8").ClearContents()

$altData = @(
    ,@("IAU","Metal","iShares Gold Trust","NYSEArca","",0.0025,0,14.52)
    ,@("DBMF","Managed Futures","iMGP DBi Managed Futures Strategy ETF","NYSEArca","Managed Futures",0.0085,0.0781,0)
    ,@("MCRO","Global Macro","IQ Hedge Macro Tracker ETF","NYSEArca","",0.0067,0.025,6.59)
    ,@("VIXM","VIX Futures","ProShares VIX Mid-Term Futures ETF","Cboe US","Trading--Miscellaneous",0.0093,0,38.36)
)

for ($i = 0; $i -lt $altData.Count; $i++) {
    $row = $i + 2
    $rec = $altData[$i]
    for ($c = 0; $c -lt 8; $c++) {
        $sheetAlt.Cells.Item($row, $c + 1).Value = $rec[$c]
    }
}

# ---------------------------------------------------------------------
# Rename the tabs last. Go through temporary names first since the
# target names are a permutation of the current names and would
# otherwise collide mid-sequence.
# ---------------------------------------------------------------------
$sheetEquity.Name = "__tmp1__"
$sheetBond.Name   = "__tmp2__"
$sheetAlt.Name    = "__tmp3__"

$sheetEquity.Name = "Equity"
$sheetBond.Name   = "Bond"
$sheetAlt.Name    = "Alternative"
